# "arreglos de la presentación" — add slide transitions.
#
# The canonical edit wraps each of these four slides' <p:transition> in an
# <mc:AlternateContent> block: a PowerPoint-2010/2013 "Choice" branch carrying
# an exotic preset (p14:ferris / p15:prstTrans drape|crush|origami) plus a
# plain <p:fade/> <mc:Fallback> that every older host (and this COM host)
# actually renders. This runtime's SlideShowTransition object only models the
# classic PpEntryEffect/Speed/Duration surface (no p14/p15 preset writer), so
# we apply the fidelity-preserving part of the edit: same slow speed, same
# per-slide duration, and the fallback's own effect (fade) — i.e. exactly
# what every one of the four <mc:Fallback> branches below already specifies.

$p = $ppt.ActivePresentation

$ppEffectFade           = [Microsoft.Office.Interop.PowerPoint.PpEntryEffect]::ppEffectFade
$ppTransitionSpeedSlow  = [Microsoft.Office.Interop.PowerPoint.PpTransitionSpeed]::ppTransitionSpeedSlow

# Slide 1 - title slide: fallback of the "ferris" (p14:ferris dir="l") choice.
$t1 = $p.Slides.Item(1).SlideShowTransition
$t1.Duration = 2
$t1.EntryEffect = $ppEffectFade
$t1.Speed = $ppTransitionSpeedSlow

# Slide 2 - "¿Qué es un singleton?": fallback of the "drape" (p15:prstTrans) choice.
$t2 = $p.Slides.Item(2).SlideShowTransition
$t2.Duration = 2
$t2.EntryEffect = $ppEffectFade
$t2.Speed = $ppTransitionSpeedSlow

# Slide 3 - "Pros y contras": fallback of the "crush" (p15:prstTrans) choice.
$t3 = $p.Slides.Item(3).SlideShowTransition
$t3.Duration = 2
$t3.EntryEffect = $ppEffectFade
$t3.Speed = $ppTransitionSpeedSlow

# Slide 4 - "Implementación en C# y Unity": fallback of the "origami" choice
# (this one runs longer in the source deck: 3.25s).
$t4 = $p.Slides.Item(4).SlideShowTransition
$t4.Duration = 3.25
$t4.EntryEffect = $ppEffectFade
$t4.Speed = $ppTransitionSpeedSlow
